$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Re-style the existing trade_date column (I2:I151): switch from the
#    date-only style (s=3 / YYYY-MM-DD) to the date-time style (s=2 /
#    YYYY-MM-DD HH:MM:SS), matching the datetime column (B).
$ws.Range("I2:I151").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 2. Append the 4 new trading-day rows (152-155). The new rows keep the
#    "old" pairing: datetime column (B) uses the datetime style, trade_date
#    column (I) uses the plain date style.
$newRows = @(
    @(2184,   46049, "NSE", 2201,   2121.1, 2130,   "SIEENE", 421850,  46049, "INE1NPP01017", "Siemens Energy India Ltd", "SIEENE", "BREEZE"),
    @(2356.9, 46050, "NSE", 2367,   2201,   2205,   "SIEENE", 1000470, 46050, "INE1NPP01017", "Siemens Energy India Ltd", "SIEENE", "BREEZE"),
    @(2481.5, 46051, "NSE", 2546,   2360.1, 2360.1, "SIEENE", 1012083, 46051, "INE1NPP01017", "Siemens Energy India Ltd", "SIEENE", "BREEZE"),
    @(2487.2, 46052, "NSE", 2519.8, 2430,   2490.3, "SIEENE", 507181,  46052, "INE1NPP01017", "Siemens Energy India Ltd", "SIEENE", "BREEZE")
)

$rowIndex = 152
foreach ($row in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]

    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]

    $ws.Cells.Item($rowIndex, 9).Value = $row[8]
    $ws.Cells.Item($rowIndex, 9).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($rowIndex, 10).Value = $row[9]
    $ws.Cells.Item($rowIndex, 11).Value = $row[10]
    $ws.Cells.Item($rowIndex, 12).Value = $row[11]
    $ws.Cells.Item($rowIndex, 13).Value = $row[12]

    $rowIndex++
}
